$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.014089
$ws.Range("H2").Value = 90.042267
$ws.Range("I2").Value = 0.5469606268302545
$ws.Range("J2").Value = 0.5469606268302545
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 1832.108205191631
$ws.Range("R2").Value = 16488.97384672468
$ws.Range("S2").Value = 0.1117776099394593
$ws.Range("T2").Value = 0.1117776099394593

$ws.Range("G3").Value = 30.014089
$ws.Range("H3").Value = 90.042267
$ws.Range("I3").Value = 0.5469606268302545
$ws.Range("J3").Value = 0.5469606268302545
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 3190.931844511474
$ws.Range("R3").Value = 28718.38660060327
$ws.Range("S3").Value = 0.194679950697506
$ws.Range("T3").Value = 0.194679950697506

$ws.Range("G4").Value = 30.014089
$ws.Range("H4").Value = 90.042267
$ws.Range("I4").Value = 0.5469606268302545
$ws.Range("J4").Value = 0.5469606268302545
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 3942.002706849097
$ws.Range("R4").Value = 35478.02436164187
$ws.Range("S4").Value = 0.2405030661932892
$ws.Range("T4").Value = 0.2405030661932892

$ws.Range("I5").Value = 0.2046507965132272
$ws.Range("J5").Value = 0.2046507965132272
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 685.5016341189894
$ws.Range("R5").Value = 6169.514707070905
$ws.Range("S5").Value = 0.0418227122471731
$ws.Range("T5").Value = 0.0418227122471731

$ws.Range("I6").Value = 0.2046507965132272
$ws.Range("J6").Value = 0.2046507965132272
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.07284145333511347
$ws.Range("T6").Value = 0.07284145333511348

$ws.Range("I7").Value = 0.2046507965132272
$ws.Range("J7").Value = 0.2046507965132272
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 1474.939793178805
$ws.Range("R7").Value = 13274.45813860925
$ws.Range("S7").Value = 0.08998663093094057
$ws.Range("T7").Value = 0.08998663093094059

$ws.Range("G8").Value = 13.63015266666667
$ws.Range("H8").Value = 40.890458
$ws.Range("I8").Value = 0.2483885766565184
$ws.Range("J8").Value = 0.2483885766565184
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 832.0064133419004
$ws.Range("R8").Value = 7488.057720077104
$ws.Range("S8").Value = 0.05076102387082106
$ws.Range("T8").Value = 0.05076102387082106

$ws.Range("G9").Value = 13.63015266666667
$ws.Range("H9").Value = 40.890458
$ws.Range("I9").Value = 0.2483885766565184
$ws.Range("J9").Value = 0.2483885766565184
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 1449.082402255143
$ws.Range("R9").Value = 13041.74162029628
$ws.Range("S9").Value = 0.08840906179581685
$ws.Range("T9").Value = 0.08840906179581687

$ws.Range("G10").Value = 13.63015266666667
$ws.Range("H10").Value = 40.890458
$ws.Range("I10").Value = 0.2483885766565184
$ws.Range("J10").Value = 0.2483885766565184
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 1790.162570210491
$ws.Range("R10").Value = 16111.46313189442
$ws.Range("S10").Value = 0.1092184909898805
$ws.Range("T10").Value = 0.1092184909898805
